# Remove StringNull enum option
#
# The "StringNull of" / "null" row (row 6 on the AllPropertyTypes sheet)
# is no longer a valid option, so delete the entire row. This shifts
# every row below it up by one, which is exactly what the target
# workbook shows (e.g. "DateTime of" moves from row 7 to row 6, the
# final "When"/"Then"/assert rows move from 16/18/19 to 15/17/18, etc.)
# and also drops the now-unused "StringNull of" / "\"null\"" shared
# strings automatically when Excel re-saves the file.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AllPropertyTypes")

# Delete row 6 ("StringNull of" / "null") - everything below shifts up.
$ws.Rows(6).Delete()

# Leave the selection where the author left it after performing the
# deletion in the UI.
$ws.Range("C19").Select()
